# Fix "Recorded By" (column G) entries: the automated "System" account
# name was accidentally listed first in the comma-separated recorder list.
# Move it to the end of the list (swap the first and last entries) for
# every row whose recorder list begins with "System, backup@backdoor.com"
# or "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $needsSwap = $val.StartsWith("System, backup@backdoor.com") -or $val.StartsWith("System, dnasr281@gmail.com")

    if ($needsSwap) {
        $parts = $val.Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        if ($trimmed.Length -ge 2) {
            $first = $trimmed[0]
            $last = $trimmed[$trimmed.Length - 1]
            $trimmed[0] = $last
            $trimmed[$trimmed.Length - 1] = $first

            $newVal = [string]::Join(", ", $trimmed)
            $cell.Value = $newVal
        }
    }
}
